$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("M15").Value = -40
$ws.Range("N15").Value = 200
$ws.Range("K14").Copy()
$ws.Range("N15").PasteSpecial(-4122)

# Row 16
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 20
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 30
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = -25
$ws.Range("L16").Value = -3.225806451612
$ws.Range("M16").Value = -9.090909090909
$ws.Range("N16").Value = -73.214285714285

# Row 17
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -69.230769230769
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 37
$ws.Range("H17").Value = -43.243243243243
$ws.Range("I17").Value = 37
$ws.Range("J17").Value = 64
$ws.Range("K17").Value = -42.1875
$ws.Range("L17").Value = -24.489795918367
$ws.Range("M17").Value = 27.586206896551
$ws.Range("N17").Value = -50.666666666666

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 20
$ws.Range("J18").Value = 12
$ws.Range("K18").Value = 66.666666666666
$ws.Range("L18").Value = 122.222222222222
$ws.Range("M18").Value = 17.647058823529
$ws.Range("N18").Value = -70.588235294117

# Row 19
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -30
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = -19.354838709677
$ws.Range("I19").Value = 49
$ws.Range("J19").Value = 54
$ws.Range("K19").Value = -9.259259259259
$ws.Range("L19").Value = -2
$ws.Range("M19").Value = 63.333333333333
$ws.Range("N19").Value = 53.125

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("J14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("F20").Value = 4
$ws.Range("H20").Value = -20
$ws.Range("J20").Value = 9
$ws.Range("K20").Value = 0
$ws.Range("N20").Value = -80.434782608695

# Row 21
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -38.709677419354
$ws.Range("F21").Value = 74
$ws.Range("G21").Value = 103
$ws.Range("H21").Value = -28.155339805825
$ws.Range("I21").Value = 148
$ws.Range("J21").Value = 183
$ws.Range("K21").Value = -19.125683060109
$ws.Range("L21").Value = 0.680272108843
$ws.Range("M21").Value = 19.354838709677
$ws.Range("N21").Value = -55.952380952380

# Row 22
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -80

# Row 23
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 9
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = -30.769230769230
$ws.Range("J23").Value = 27
$ws.Range("K23").Value = -3.703703703703
$ws.Range("L23").Value = 52.941176470588

# Row 24
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = -40.909090909090
$ws.Range("F24").Value = 63
$ws.Range("G24").Value = 55
$ws.Range("H24").Value = 14.545454545454
$ws.Range("I24").Value = 107
$ws.Range("J24").Value = 110
$ws.Range("K24").Value = -2.727272727272
$ws.Range("L24").Value = -48.309178743961
$ws.Range("M24").Value = -28.666666666666

# Row 25
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -25
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = 43.75
$ws.Range("I25").Value = 32
$ws.Range("J25").Value = 32
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = -70.093457943925

# Row 26
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -36.363636363636
$ws.Range("F26").Value = 36
$ws.Range("H26").Value = -5.263157894736
$ws.Range("I26").Value = 92
$ws.Range("J26").Value = 77
$ws.Range("K26").Value = 19.480519480519
$ws.Range("L26").Value = 27.777777777777
$ws.Range("M26").Value = 39.393939393939

# Row 27
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

# Row 28
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = -40
$ws.Range("L28").Value = 0

# Row 29
$ws.Range("L29").Value = -60
$ws.Range("N29").Value = -50

# Row 30
$ws.Range("L30").Value = -60
$ws.Range("N30").Value = -50

# Row 33
$ws.Range("F33").NumberFormat = "@"
$ws.Range("F33").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F33").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Header text updates (Volume number, date range)
$volOld = $ws.Range("A8").Value()
$volNew = $volOld.Substring(0, $volOld.Length - 1) + "8"
$ws.Range("A8").Value = $volNew

$dateOld = $ws.Range("C9").Value()
$dateNew = $dateOld.Replace("2/10/2025", "2/17/2025").Replace("2/16/2025", "2/23/2025")
$ws.Range("C9").Value = $dateNew
